$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "magapoke_2025-12-31"

$ws.Range("A1").Value = "rank"
$ws.Range("B1").Value = "title"
$headerRange = $ws.Range("A1:B1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$titles = @(
    'WIND BREAKER',
    'ブルーロック',
    'みいちゃんと山田さん',
    '東京卍リベンジャーズ',
    'ギルティサークル',
    'ベイビーステップ',
    '信じていた仲間達にダンジョン奥地で殺されかけたがギフト『無限ガチャ』でレベル9999の仲間達を手に入れて元パーティーメンバーと世界に復讐＆『ざまぁ！』します！',
    '島耕作',
    '十字架のろくにん',
    'ガチアクタ',
    'イレギュラーズ',
    '黄昏町プリズナーズ',
    '薫る花は凛と咲く',
    '愛妻の裏アカ',
    'ハードワーカー中田',
    '南海トラフ巨大地震',
    '黒猫と魔女の教室',
    '君が僕らを悪魔と呼んだ頃',
    '転生したら第七王子だったので、気ままに魔術を極めます',
    '魔術ギルド総帥～生まれ変わって今更やり直す2度目の学院生活～',
    'となりの黒川さん',
    '魔女と傭兵',
    '蒼く染めろ',
    '転生貴族、鑑定スキルで成り上がる～弱小領地を受け継いだので、優秀な人材を増やしていたら、最強領地になってた～',
    '異世界ウォーキング',
    'K-9~警視庁公安部公安第9課異能対策係~',
    '味方が弱すぎて補助魔法に徹していた宮廷魔法師、追放されて最強を目指す',
    'アルキメデスの大戦',
    'Fate/Grand Order -Epic of Remnant- 英霊剣豪七番勝負',
    'FAIRY TAIL 100 YEARS QUEST',
    '食糧人類-Starving Anonymous-',
    '幼馴染とはラブコメにならない',
    'グラぱらっ！',
    'ハナバス　苔石花江のバスケ論',
    'せいぶつ部の田辺くん',
    'ひゃくえむ。',
    'ナキナギ',
    'さわらないで小手指くん',
    'ともだちづくり',
    '田んぼで拾った女騎士、田舎で俺の嫁だと思われている',
    'おやすみ ふみさん',
    '追放された転生王子、『自動製作《オートクラフト》』スキルで領地を爆速で開拓し最強の村を作ってしまう〜最強クラフトスキルで始める、楽々領地開拓スローライフ〜',
    '降り積もれ孤独な死よ',
    '追放されなかった男　～二度目の人生は土下座から始まりました～',
    'お母さん冒険者、ログインボーナスでスキル【主婦】に目覚めました。週一貰えるチラシで冒険者生活頑張ります！',
    'ドラハチ',
    '皇女転生　～伝説の大魔導士（♂）、姫騎士となりて伝説の令嬢騎士団を作り無双する～',
    'ダメスキル【自動機能】が覚醒しました～あれ、ギルドのスカウトの皆さん、俺を「いらない」って言ってませんでした？～',
    'この世界がいずれ滅ぶことを、俺だけが知っている～モンスターが現れた世界で、死に戻りレベルアップ～',
    '春くらり',
    '時々ボソッとロシア語でデレる隣のアーリャさん',
    'デッドアカウント',
    'Aランクパーティを離脱した俺は、元教え子たちと迷宮深部を目指す。',
    'ジュミドロ',
    '限界集落を脱村した錬金術士、都会で"最強"なのがバレまくる。～老害どもにはいい加減愛想が尽きました～',
    'ストーカー行為がバレて人生終了男',
    'アオバノバスケ',
    '異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～',
    '辺境の薬師、都でSランク冒険者となる～英雄村の少年がチート薬で無自覚無双〜',
    '不遇職『鍛冶師』だけど最強です ～気づけば何でも作れるようになっていた男ののんびりスローライフ～',
    '阿武ノーマル',
    'ルックスＹを選んでしまいました 〜やり込んでいるゲームに転生したはずなのに、未実装のガチャで攻略をすることになった件〜',
    '念願の悪役令嬢（ラスボス）の身体を手に入れたぞ！',
    '魁の花巫女',
    'デスティニーラバーズ',
    '屋根の下のアルテミス',
    'いじめるヤバイ奴',
    '普通の本はありません！',
    'はっちぽっちぱんち',
    '東京卍リベンジャーズ～場地圭介からの手紙～',
    'MYS',
    'リスナーに騙されてダンジョンの最下層から脱出RTAすることになった',
    'Destiny Unchain Online 〜吸血鬼少女となって、やがて『赤の魔王』と呼ばれるようになりました〜',
    '可愛いだけじゃない式守さん',
    'なれの果ての僕ら',
    '我間乱 ―修羅―',
    '劣等人の魔剣使い　スキルボードを駆使して最強に至る',
    '不遇職【鑑定士】が実は最強だった～奈落で鍛えた最強の【神眼】で無双する～',
    'それがメイドのカンナです',
    '君が監督！',
    '死ぬほど君の処女が欲しい',
    '母という呪縛 娘という牢獄',
    'ヒロインは絶望しました。',
    '中華一番！極',
    '卒業アルバムの彼女たち',
    'イジらないで、長瀞さん',
    '捨てられた転生賢者～魔物の森で最強の大魔帝国を作り上げる～',
    'GALAXIAS',
    'メイドの岸さん',
    'ウイニング パス',
    '英雄と魔女の転生ラブコメ',
    '五輪の女神さま 〜なでしこ寮のメダルごはん〜',
    '復讐の教科書',
    '絶望の楽園',
    '金田一少年の事件簿外伝 犯人たちの事件簿',
    '四刀流の最強配信者～やり込んだVRゲームの設定が現実世界に反映されたので、廃止予定だった戦闘職で無双します～',
    'JK Biker',
    'インフェクション',
    '鳴るさんだぁ',
    'ハプスブルク家の華麗なる受難'
)

for ($i = 0; $i -lt $titles.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $titles[$i]
}

$wb.Worksheets.Item(1).Activate()

Write-Output ("Sheet count: " + $wb.Worksheets.Count)
Write-Output ("New sheet name: " + $ws.Name)
